$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Saisie")

# Update the journée/week number (P42 -> P47), correction de Samuel B.
$ws.Range("B2").Value = 47

# The K5:K13 shared formulas read their F/G/H/I operands from the cached
# "spill" values produced by the F4 array formula (FILTER). Touch each
# shared-formula cell so its cached result is refreshed against the new
# spilled values (re-assigning the same formula text forces a recompute
# without altering the on-disk shared-formula grouping).
foreach ($addr in @("K5","K6","K7","K8","K9","K10","K11","K12","K13")) {
    $cell = $ws.Range($addr)
    $cell.Formula = $cell.Formula
}

$excel.CalculateFullRebuild()
$excel.Calculate()

# Reflect the active selection left by the author after the edit
$ws.Activate()
$ws.Range("B3").Select()
